$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.381.86"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "3.662.56"
$ws.Range("E3").Value = "  +2.07%  "

$ws.Range("E4").Value = "  -0.03%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "241.86"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.90"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +14.92%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "659.46"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.76%  "

$ws.Range("E8").Value = "  +2.78%  "

$ws.Range("E9").Value = "  +2.84%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"

$ws.Range("D11").Value = "3.661.42"
$ws.Range("E11").Value = "  +2.05%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "44.94"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.81%  "

$ws.Range("E13").Value = "  +0.32%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.69"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.85%  "

$ws.Range("D15").Value = "4.342.98"
$ws.Range("E15").Value = "  +2.11%  "

$ws.Range("E16").Value = "  +5.26%  "

$ws.Range("D17").Value = "96.155.56"
$ws.Range("E17").Value = "  -0.71%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "8.89"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +14.53%  "

$ws.Range("D19").Value = "3.658.94"
$ws.Range("E19").Value = "  +2.22%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.73"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.09%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "18.27"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.534"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "521.51"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.36%  "

$ws.Range("E25").Value = "  +1.27%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.90"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "102.37"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.48%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "12.99"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.38%  "

$ws.Range("E29").Value = "  +8.34%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "12.29"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +6.64%  "

$ws.Range("E31").Value = "  -0.90%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.185"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("E34").Value = "  +11.43%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "33.13"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +5.18%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.82%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.587"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.72%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "625.66"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.23%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "46.12"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +38.52%  "

$ws.Range("E40").Value = "  -2.57%  "

$ws.Range("E41").Value = "  +5.13%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.959"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +4.83%  "

$ws.Range("E43").Value = "  +5.46%  "

$ws.Range("E44").Value = "  -0.02%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "6.27"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +7.66%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0453"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.30%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.427"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +18.62%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.29"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.81%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "23.60"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "8.52"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.59%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "3.59"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.68%  "
